$wb = $excel.ActiveWorkbook

# --- helper: write a value to a cell while forcing TEXT storage, even when
# the text looks like a plain number (e.g. "-2.9"). Excel normally
# auto-converts a numeric-looking string typed into ".Value" into a real
# number; to avoid that we stage the text in a scratch cell that has been
# explicitly formatted as Text, copy it, and paste-special "values only"
# into the destination. Because the destination cell's own style/format is
# never touched, it keeps using the workbook's default style (no stray
# number-format / style entries end up attached to the edited cells).
function Set-TextValue($Worksheet, $Cell, $Text) {
    $scratch = $Worksheet.Range("ZZ1000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy()
    $Worksheet.Range($Cell).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# Restricciones_del_lider
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
$ws.Range("A2").Value = "1.9 - x"
Set-TextValue $ws "B2" "-2.9"
Set-TextValue $ws "D2" "0.83"
$ws.Range("A3").Value = "-1.9 + x"
Set-TextValue $ws "B3" "0.8999999999999999"
Set-TextValue $ws "D3" "0.08"

# ---------------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "-0.30000000000000004 + 2.0y"
Set-TextValue $ws "B2" "-0.7"
Set-TextValue $ws "D2" "0.01"
Set-TextValue $ws "E2" "1.3"
Set-TextValue $ws "F2" "3.9000000000000004"
Set-TextValue $ws "A3" "0"
Set-TextValue $ws "B3" "-1"
Set-TextValue $ws "D3" "0.97"
Set-TextValue $ws "E3" "2.4"
Set-TextValue $ws "F3" "3.2"

# ---------------------------------------------------------------------
# Punto_modificado
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws "A2" "1.9"
Set-TextValue $ws "B2" "0.15"

# ---------------------------------------------------------------------
# Vector_bf
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $ws "A2" "0.5432499999999999"

# ---------------------------------------------------------------------
# Vector_BF
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $ws "A2" "2.185"
Set-TextValue $ws "A3" "-0.3532499999999996"

# ---------------------------------------------------------------------
# Vector_Alpha (this one is a genuine numeric cell, not text)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 1.9500000000000002
